$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-15 Friday" "2024-11-16 Saturday"

Replace-Text "302×2=" "291×9="
Replace-Text "599×9=" "247×8="
Replace-Text "323×5=" "714×2="
Replace-Text "851×7=" "977×3="
Replace-Text "190×7=" "218×6="

Replace-Text "239×2=" "612×5="
Replace-Text "917×8=" "683×2="
Replace-Text "678×4=" "560×4="
Replace-Text "499×6=" "767×9="
Replace-Text "249×5=" "838×6="

Replace-Text "784×6=" "865×8="
Replace-Text "224×2=" "113×8="
Replace-Text "903×3=" "424×7="
Replace-Text "269×7=" "670×5="
Replace-Text "974×6=" "244×2="

Replace-Text "189×5=" "715×8="
Replace-Text "470×9=" "909×8="
Replace-Text "194×4=" "571×3="
Replace-Text "970×9=" "746×6="
Replace-Text "186×8=" "546×9="

Replace-Text "581×3=" "169×3="
Replace-Text "666×6=" "153×3="
Replace-Text "391×5=" "660×2="
Replace-Text "716×3=" "576×4="
Replace-Text "717×6=" "478×3="

Write-Host "Done"
